$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 394.14285
$ws.Range("I9").Value = 384.83334
$ws.Range("J9").Value = 450
$ws.Range("K9").Value = 384.83334
$ws.Range("L9").Value = 450
$ws.Range("M9").Value = -215.83334
$ws.Range("N9").Value = -788

$ws.Range("H18").Value = 1371.125
$ws.Range("I18").Value = 852.7143
$ws.Range("K18").Value = 852.7143
$ws.Range("M18").Value = -568.7143

$ws.Range("H32").Value = 3111.818
$ws.Range("J32").Value = 3219
$ws.Range("L32").Value = 3219
$ws.Range("N32").Value = -3871

$ws.Range("H33").Value = 249.72728
$ws.Range("J33").Value = 346.33334
$ws.Range("L33").Value = 346.33334
$ws.Range("N33").Value = -804.33334

$ws.Range("H70").Value = 15000
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 15000
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws.Range("H100").Value = 3001
$ws.Range("J100").Value = 3001
$ws.Range("L100").Value = 3001
$ws.Range("N100").Value = -4083

$ws.Range("H112").Value = 3563.7693
$ws.Range("I112").Value = 2796.5
$ws.Range("J112").Value = 3703.2727
$ws.Range("K112").Value = 8389.5
$ws.Range("L112").Value = 11109.8181
$ws.Range("M112").Value = -7281.5
$ws.Range("N112").Value = -13325.8181

$ws.Range("H115").Value = 2999
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 2999
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 8997
$ws.Range("M115").ClearContents()
$ws.Range("N115").Value = -12131

$ws.Range("H116").Value = 3816.36
$ws.Range("I116").Value = 3039.1538
$ws.Range("J116").Value = 4658.3335
$ws.Range("K116").Value = 3039.1538
$ws.Range("L116").Value = 4658.3335
$ws.Range("M116").Value = 402.8462
$ws.Range("N116").Value = -11542.3335

$ws.Range("H137").Value = 1311.1666
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()

$ws.Range("H138").Value = 3857.9788
$ws.Range("J138").Value = 4499.25
$ws.Range("L138").Value = 13497.75
$ws.Range("N138").Value = -23777.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2335435.2
$ws.Range("I32").Value = 2189667.5
$ws.Range("K32").Value = 2189667.5
$ws.Range("M32").Value = -2189380.5

$ws.Range("H61").Value = 10634.056
$ws.Range("I61").Value = 10789
$ws.Range("J61").Value = 8000
$ws.Range("K61").Value = 10789
$ws.Range("L61").Value = 8000
$ws.Range("M61").Value = -10577
$ws.Range("N61").Value = -8424

$ws.Range("H97").Value = 1047.7142
$ws.Range("I97").Value = 1132.8
$ws.Range("K97").Value = 1132.8
$ws.Range("M97").Value = -636.8

$ws.Range("H132").Value = 2199
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws.Range("H136").Value = 10634.056
$ws.Range("I136").Value = 10789
$ws.Range("J136").Value = 8000
$ws.Range("K136").Value = 32367
$ws.Range("L136").Value = 24000
$ws.Range("M136").Value = -29817
$ws.Range("N136").Value = -29100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2828
$ws.Range("I94").Value = 2556.875
$ws.Range("J94").Value = 4997
$ws.Range("K94").Value = 2556.875
$ws.Range("L94").Value = 4997
$ws.Range("M94").Value = -2105.875
$ws.Range("N94").Value = -5899

$ws.Range("H107").Value = 5766.3335
$ws.Range("I107").Value = 5766.3335
$ws.Range("K107").Value = 5766.3335
$ws.Range("M107").Value = -3846.3335

$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1874.5
$ws.Range("I22").Value = 1832.6666
$ws.Range("K22").Value = 1832.6666
$ws.Range("M22").Value = -1482.6666

$ws.Range("H31").Value = 1675.9474
$ws.Range("J31").Value = 1711.7858
$ws.Range("L31").Value = 1711.7858
$ws.Range("N31").Value = -2301.7858

$ws.Range("H34").Value = 1675.9474
$ws.Range("J34").Value = 1711.7858
$ws.Range("L34").Value = 1711.7858
$ws.Range("N34").Value = -2115.7858

$ws.Range("H51").Value = 50000
$ws.Range("J51").Value = 70000
$ws.Range("L51").Value = 70000
$ws.Range("N51").Value = -71472

$ws.Range("H58").Value = 5066.2
$ws.Range("I58").Value = 3259
$ws.Range("K58").Value = 3259
$ws.Range("M58").Value = -3056

$ws.Range("H61").Value = 50000
$ws.Range("J61").Value = 70000
$ws.Range("L61").Value = 70000
$ws.Range("N61").Value = -70696

$ws.Range("H107").Value = 1544.8334
$ws.Range("I107").Value = 1155.5
$ws.Range("J107").Value = 1739.5
$ws.Range("K107").Value = 1155.5
$ws.Range("L107").Value = 1739.5
$ws.Range("M107").Value = 764.5
$ws.Range("N107").Value = -5579.5

$ws.Range("H132").Value = 8406.182000000001
$ws.Range("I132").Value = 8441.777
$ws.Range("J132").Value = 8246
$ws.Range("K132").Value = 25325.331
$ws.Range("L132").Value = 24738
$ws.Range("M132").Value = -22795.331
$ws.Range("N132").Value = -29798

$ws.Range("H136").Value = 5066.2
$ws.Range("I136").Value = 3259
$ws.Range("K136").Value = 9777
$ws.Range("M136").Value = -7227

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 125.25
$ws.Range("I2").Value = 24.75
$ws.Range("K2").Value = 148.5
$ws.Range("M2").Value = -35.5

$ws.Range("H4").Value = 8637645
$ws.Range("I4").Value = 10002738
$ws.Range("K4").Value = 30008214
$ws.Range("M4").Value = -30008102

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 88.29412000000001
$ws.Range("I2").Value = 76.59999999999999
$ws.Range("J2").Value = 176
$ws.Range("K2").Value = 76.59999999999999
$ws.Range("L2").Value = 176
$ws.Range("M2").Value = 36.40000000000001
$ws.Range("N2").Value = -402

$ws.Range("H132").Value = 3942.2727
$ws.Range("I132").Value = 3942.2727
$ws.Range("K132").Value = 11826.8181
$ws.Range("M132").Value = -9296.8181

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3169.8
$ws.Range("I22").Value = 2712.5
$ws.Range("J22").Value = 4999
$ws.Range("K22").Value = 2712.5
$ws.Range("L22").Value = 4999
$ws.Range("M22").Value = -2417.5
$ws.Range("N22").Value = -5589

$ws.Range("H27").Value = 3169.8
$ws.Range("I27").Value = 2712.5
$ws.Range("J27").Value = 4999
$ws.Range("K27").Value = 2712.5
$ws.Range("L27").Value = 4999
$ws.Range("M27").Value = -2605.5
$ws.Range("N27").Value = -5213

$ws.Range("H93").Value = 967.6667
$ws.Range("I93").Value = 2003
$ws.Range("K93").Value = 2003
$ws.Range("M93").Value = -755

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 455.75
$ws.Range("J113").Value = 499.66666
$ws.Range("L113").Value = 1498.99998
$ws.Range("N113").Value = -5838.999980000001
